$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend the table (rows 256-269) with new daily records, reusing the
# formatting of the last existing row (row 255) for the date column style.
$ws.Range("A255:D255").Copy()
$ws.Range("A256:D269").PasteSpecial(-4122)

$dates    = @(44330,44331,44332,44333,44334,44335,44336,44337,44338,44339,44340,44341,44342,44343)
$newpos   = @(1,0,0,0,0,0,1,1,1,1,1,0,0,0)
$somma7gg = @(3,3,3,3,3,2,2,2,3,4,5,5,5,4)
$per100k  = @(80.29978586723769,80.29978586723769,80.29978586723769,80.29978586723769,80.29978586723769,53.53319057815846,53.53319057815846,53.53319057815846,80.29978586723769,107.0663811563169,133.8329764453961,133.8329764453961,133.8329764453961,107.0663811563169)

for ($i = 0; $i -lt $dates.Count; $i++) {
    $r = 256 + $i
    $ws.Cells.Item($r, 1).Value = $dates[$i]
    $ws.Cells.Item($r, 2).Value = $newpos[$i]
    $ws.Cells.Item($r, 3).Value = $somma7gg[$i]
    $ws.Cells.Item($r, 4).Value = $per100k[$i]
}
